# The deck currently ships two theme parts:
#   ppt/theme/theme1.xml -> "Integral" (Red Violet colours) - used by the slide master
#   ppt/theme/theme2.xml -> "Office Theme" (Office colours) - used by the notes master
#
# The target edit swaps their content so the slide master ends up using the
# standard Office Theme colour palette (what used to live in theme2.xml).
#
# PowerPoint's ColorScheme/ThemeColorScheme object model addresses theme
# colours positionally:
#   1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
#   9=accent5 10=accent6 11=hlink 12=folHlink
# and .RGB takes a standard VBA/COM color value (0xBBGGRR), so we push the
# Office Theme's RGB values (expressed as 0xBBGGRR longs) onto the slide
# master's theme color scheme.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Colors(1).RGB  = 0            # dk1      000000
$cs.Colors(2).RGB  = 16777215     # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388       # dk2      44546A
$cs.Colors(4).RGB  = 15132391      # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939      # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501       # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845      # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407         # accent4  FFC000
$cs.Colors(9).RGB  = 12874308      # accent5  4472C4
$cs.Colors(10).RGB = 4697456       # accent6  70AD47
$cs.Colors(11).RGB = 12673797      # hlink    0563C1
$cs.Colors(12).RGB = 7491477       # folHlink 954F72
